$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.842.03"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.18"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.72"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5343"
$ws.Range("E7").Value = "  +2.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3753"
$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07178"
$ws.Range("E9").Value = "  -1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.62"
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8858"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08114"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.60"
$ws.Range("E13").Value = "  +1.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.03"
$ws.Range("E14").Value = "  -2.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.271"
$ws.Range("E15").Value = "  -1.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.71"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008534"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.892.76"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.972"
$ws.Range("E21").Value = "  -3.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.387"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.12"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.245"
$ws.Range("E25").Value = "  -3.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.733"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.40"
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.744"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.570"
$ws.Range("E30").Value = "  -6.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09136"
$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7972"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04984"
$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.996"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.168"
$ws.Range("E35").Value = "  -4.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5901"
$ws.Range("E36").Value = "  +3.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.606"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.151"
$ws.Range("E38").Value = "  -6.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01948"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.652"
$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.904"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.65"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5041"
$ws.Range("E44").Value = "  +3.08%  "

$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.973"
$ws.Range("E47").Value = "  -2.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.620"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.64"
$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("E51").Value = "  -2.73%  "
